$wb = $excel.ActiveWorkbook

# Build the new sheet by duplicating an existing sheet that already has the
# exact same "Index / <header>" layout (2 data rows + header, same styles,
# same page setup), then adjust its name/values. This keeps formatting
# (header fill/bold, quote-prefixed index cells, page setup, etc.) consistent
# with the rest of the workbook.
$template = $wb.Worksheets.Item("Live_Tv_Tab_Live_And_Catch")
$afterSheet = $wb.Worksheets.Item("Continue_Watching")
$template.Copy([System.Reflection.Missing]::Value, $afterSheet) | Out-Null

$newSheet = $wb.Worksheets.Item($afterSheet.Index + 1)
$newSheet.Name = "List_Tab_Live_Event"

# Column A (Index / 0 / 1) already matches what we need from the template.
# Update column B with the new tab values.
$newSheet.Range("B1").Value = "Tabs"
$newSheet.Range("B2").Value = "Live Event"
$newSheet.Range("B3").Value = "Missed Event"

# Normalize the selection and make this newly inserted sheet the active tab
$newSheet.Range("A1").Select() | Out-Null
$newSheet.Activate()
